# Slide 5 ("ch9.pptx" slide index 5) contains a text box (shape #4,
# "Text Box 12") whose sole run reads "Notes:". The edit appends a
# trailing space so it reads "Notes: ".
#
# The shape has <a:spAutoFit/> in its body properties, so simply
# reassigning TextRange.Text causes the COM host to re-layout the box
# and recompute its height (Shape.Height) to a new autofit value. The
# canonical edit captured in the diff does NOT touch the shape's
# size/position (<a:ext cx=".." cy=".."/> is unchanged), so we restore
# the original height immediately afterwards to avoid that incidental
# side effect.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(5)
$sh = $s.Shapes.Item(4)

$originalHeight = $sh.Height

$sh.TextFrame.TextRange.Text = "Notes: "

$sh.Height = $originalHeight
